$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1: sheet structure.
# The workbook currently ends in "...,2021-Q4,总计" (总计 = sheetId 6).
# Target: "...,2021-Q4,2022-Q1,总计" where 2022-Q1 keeps sheetId 6 and
# 总计 gets a fresh sheetId 7 (matches how Excel assigns ids: renaming
# preserves the id, Add() mints a new one).
#
# Since the existing "总计" sheet is reused (renamed) to become
# "2022-Q1" and its old grid gets wiped, its original 5 summary rows
# are hard-coded below and replayed into the brand-new "总计" sheet.
# ------------------------------------------------------------------

# 1a. Rename the existing "总计" sheet in place -> "2022-Q1". It keeps
#     its original sheetId and tab position (last). We repopulate its
#     contents with the fund-holding detail further down.
$wb.Worksheets.Item("总计").Name = "2022-Q1"

# 1b. Add a brand-new sheet (Add() always inserts at position 1 in this
#     runtime) then move it so it lands right after "2022-Q1", i.e. at
#     the very end, and rename it to "总计".
#     NOTE: sheet-collection handles captured in a variable resolve by
#     POSITION, not identity - any Add/Move/Delete invalidates earlier
#     captures. So every lookup below is re-fetched fresh right before use.
$wb.Worksheets.Add() | Out-Null
$wb.Worksheets.Item(1).Move($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$wb.Worksheets.Item("Sheet1").Name = "总计"

# ------------------------------------------------------------------
# Step 2: populate the "2022-Q1" sheet with fund-holding detail rows
# (same shape as the other quarterly sheets, e.g. "2021-Q4").
# ------------------------------------------------------------------
$wsQ = $wb.Worksheets.Item("2022-Q1")

# Wipe the old (总计-shaped) content + formatting first.
$wsQ.Cells.Clear()

# Borrow the exact header / index-column formatting (bold, border,
# center/top alignment -> style index 2 in the sibling sheets) from the
# "2021-Q4" sheet, which already has the identical 7-column layout.
$wb.Worksheets.Item("2021-Q4").Range("B1:H1").Copy()
$wb.Worksheets.Item("2022-Q1").Range("B1:H1").PasteSpecial(-4122)
$wb.Worksheets.Item("2021-Q4").Range("A2").Copy()
$wb.Worksheets.Item("2022-Q1").Range("A2:A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row.
$wsQ.Range("B1").Value = "基金代码"
$wsQ.Range("C1").Value = "基金名称"
$wsQ.Range("D1").Value = "基金规模"
$wsQ.Range("E1").Value = "股票总仓位"
$wsQ.Range("F1").Value = "仓位占比"
$wsQ.Range("G1").Value = "持有市值(亿元)"
$wsQ.Range("H1").Value = "仓位排名"

# Force columns B:G to store plain text (fund codes must keep leading
# zeroes, and the size/position/value columns are text in the source
# data too) - otherwise Excel auto-coerces numeric-looking strings.
$wsQ.Range("B2:G12").NumberFormat = "@"

# A-column running index (0-based) and data rows.

# Row 2: 012930 中庚价值先锋股票
$wsQ.Range("A2").Value = 0
$wsQ.Range("B2").Value = "012930"
$wsQ.Range("C2").Value = "中庚价值先锋股票"
$wsQ.Range("D2").Value = "54.59"
$wsQ.Range("E2").Value = "94.46"
$wsQ.Range("F2").Value = "5.27"
$wsQ.Range("G2").Value = "2.8769"
$wsQ.Range("H2").Value = 6

# Row 3: 009230 鹏华安和混合A
$wsQ.Range("A3").Value = 1
$wsQ.Range("B3").Value = "009230"
$wsQ.Range("C3").Value = "鹏华安和混合A"
$wsQ.Range("D3").Value = "14.02"
$wsQ.Range("E3").Value = "34.45"
$wsQ.Range("F3").Value = "1.36"
$wsQ.Range("G3").Value = "0.1907"
$wsQ.Range("H3").Value = 5

# Row 4: 009667 鹏华安庆混合A
$wsQ.Range("A4").Value = 2
$wsQ.Range("B4").Value = "009667"
$wsQ.Range("C4").Value = "鹏华安庆混合A"
$wsQ.Range("D4").Value = "11.22"
$wsQ.Range("E4").Value = "38.92"
$wsQ.Range("F4").Value = "1.52"
$wsQ.Range("G4").Value = "0.1705"
$wsQ.Range("H4").Value = 5

# Row 5: 013393 信达澳银价值精选混合A
$wsQ.Range("A5").Value = 3
$wsQ.Range("B5").Value = "013393"
$wsQ.Range("C5").Value = "信达澳银价值精选混合A"
$wsQ.Range("D5").Value = "3.61"
$wsQ.Range("E5").Value = "81.31"
$wsQ.Range("F5").Value = "3.23"
$wsQ.Range("G5").Value = "0.1166"
$wsQ.Range("H5").Value = 4

# Row 6: 005416 鹏华尊惠18个月定期开放混合A
$wsQ.Range("A6").Value = 4
$wsQ.Range("B6").Value = "005416"
$wsQ.Range("C6").Value = "鹏华尊惠18个月定期开放混合A"
$wsQ.Range("D6").Value = "7.95"
$wsQ.Range("E6").Value = "37.81"
$wsQ.Range("F6").Value = "1.24"
$wsQ.Range("G6").Value = "0.0986"
$wsQ.Range("H6").Value = 7

# Row 7: 009231 鹏华安和混合C
$wsQ.Range("A7").Value = 5
$wsQ.Range("B7").Value = "009231"
$wsQ.Range("C7").Value = "鹏华安和混合C"
$wsQ.Range("D7").Value = "5.33"
$wsQ.Range("E7").Value = "34.45"
$wsQ.Range("F7").Value = "1.36"
$wsQ.Range("G7").Value = "0.0725"
$wsQ.Range("H7").Value = 5

# Row 8: 003165 鹏华弘嘉灵活配置混合A
$wsQ.Range("A8").Value = 6
$wsQ.Range("B8").Value = "003165"
$wsQ.Range("C8").Value = "鹏华弘嘉灵活配置混合A"
$wsQ.Range("D8").Value = "1.53"
$wsQ.Range("E8").Value = "93.95"
$wsQ.Range("F8").Value = "3.74"
$wsQ.Range("G8").Value = "0.0572"
$wsQ.Range("H8").Value = 2

# Row 9: 009668 鹏华安庆混合C
$wsQ.Range("A9").Value = 7
$wsQ.Range("B9").Value = "009668"
$wsQ.Range("C9").Value = "鹏华安庆混合C"
$wsQ.Range("D9").Value = "2.36"
$wsQ.Range("E9").Value = "38.92"
$wsQ.Range("F9").Value = "1.52"
$wsQ.Range("G9").Value = "0.0359"
$wsQ.Range("H9").Value = 5

# Row 10: 003166 鹏华弘嘉灵活配置混合C
$wsQ.Range("A10").Value = 8
$wsQ.Range("B10").Value = "003166"
$wsQ.Range("C10").Value = "鹏华弘嘉灵活配置混合C"
$wsQ.Range("D10").Value = "0.56"
$wsQ.Range("E10").Value = "93.95"
$wsQ.Range("F10").Value = "3.74"
$wsQ.Range("G10").Value = "0.0209"
$wsQ.Range("H10").Value = 2

# Row 11: 013394 信达澳银价值精选混合C
$wsQ.Range("A11").Value = 9
$wsQ.Range("B11").Value = "013394"
$wsQ.Range("C11").Value = "信达澳银价值精选混合C"
$wsQ.Range("D11").Value = "0.37"
$wsQ.Range("E11").Value = "81.31"
$wsQ.Range("F11").Value = "3.23"
$wsQ.Range("G11").Value = "0.0120"
$wsQ.Range("H11").Value = 4

# Row 12: 005417 鹏华尊惠18个月定期开放混合C
$wsQ.Range("A12").Value = 10
$wsQ.Range("B12").Value = "005417"
$wsQ.Range("C12").Value = "鹏华尊惠18个月定期开放混合C"
$wsQ.Range("D12").Value = "0.56"
$wsQ.Range("E12").Value = "37.81"
$wsQ.Range("F12").Value = "1.24"
$wsQ.Range("G12").Value = "0.0069"
$wsQ.Range("H12").Value = 7

# ------------------------------------------------------------------
# Step 3: populate the brand-new "总计" sheet: header, the new
# "2022-Q1" row on top, then the 5 original summary rows (shifted
# down one row, A-index renumbered 1..5).
# ------------------------------------------------------------------
$wsT = $wb.Worksheets.Item("总计")

# Borrow formatting from the sheet that used to be "总计" (now
# "2022-Q1" holds a different layout) - instead grab it from any of
# the other quarter sheets, whose B1:D1 header + A-column share the
# very same style index (2) used by the summary sheet.
$wb.Worksheets.Item("2021-Q4").Range("B1:D1").Copy()
$wb.Worksheets.Item("总计").Range("B1:D1").PasteSpecial(-4122)
$wb.Worksheets.Item("2021-Q4").Range("A2").Copy()
$wb.Worksheets.Item("总计").Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsT.Range("B1").Value = "日期"
$wsT.Range("C1").Value = "持有数量(只)"
$wsT.Range("D1").Value = "持有市值(亿元)"

# New row: the 2022-Q1 summary.
$wsT.Range("A2").Value = 0
$wsT.Range("B2").Value = "2022-Q1"
$wsT.Range("C2").Value = 11
$wsT.Range("D2").Value = 3.66

# Row 3: 2021-Q4 (shifted down from the original row 2)
$wsT.Range("A3").Value = 1
$wsT.Range("B3").Value = "2021-Q4"
$wsT.Range("C3").Value = 20
$wsT.Range("D3").Value = 5.22

# Row 4: 2021-Q3 (shifted down from the original row 3)
$wsT.Range("A4").Value = 2
$wsT.Range("B4").Value = "2021-Q3"
$wsT.Range("C4").Value = 6
$wsT.Range("D4").Value = 1.08

# Row 5: 2021-Q2 (shifted down from the original row 4)
$wsT.Range("A5").Value = 3
$wsT.Range("B5").Value = "2021-Q2"
$wsT.Range("C5").Value = 11
$wsT.Range("D5").Value = 1.45

# Row 6: 2021-Q1 (shifted down from the original row 5)
$wsT.Range("A6").Value = 4
$wsT.Range("B6").Value = "2021-Q1"
$wsT.Range("C6").Value = 6
$wsT.Range("D6").Value = 1.42

# Row 7: 2020-Q4 (shifted down from the original row 6)
$wsT.Range("A7").Value = 5
$wsT.Range("B7").Value = "2020-Q4"
$wsT.Range("C7").Value = 7
$wsT.Range("D7").Value = 3.87

